$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Date column (BF) stays text, not auto-converted to a date serial
$ws.Range("BF2:BF31").NumberFormat = "@"

$ws.Range("BF2").Value = "2009-06-21"
$ws.Range("D3").Value = 82
$ws.Range("E3").Value = 62
$ws.Range("G3").Value = 0.756
$ws.Range("I3").Value = 37.5
$ws.Range("J3").Value = 77.2
$ws.Range("M3").Value = 16.5
$ws.Range("N3").Value = 0.397
$ws.Range("O3").Value = 19.4
$ws.Range("P3").Value = 25.3
$ws.Range("Q3").Value = 0.765
$ws.Range("S3").Value = 31.5
$ws.Range("T3").Value = 42.1
$ws.Range("U3").Value = 22.7
$ws.Range("W3").Value = 7.6
$ws.Range("X3").Value = 4.7
$ws.Range("Y3").Value = 4.7
$ws.Range("AA3").Value = 22.2
$ws.Range("AB3").Value = 100.9
$ws.Range("AD3").Value = 1
$ws.Range("AP3").Value = 11
$ws.Range("AQ3").Value = 18
$ws.Range("AS3").Value = 7
$ws.Range("AT3").Value = 9
$ws.Range("AV3").Value = 28
$ws.Range("AW3").Value = 7
$ws.Range("AX3").Value = 17
$ws.Range("AZ3").Value = 28
$ws.Range("BA3").Value = 7
$ws.Range("BB3").Value = 11
$ws.Range("BC3").Value = 3
$ws.Range("BF3").Value = "2009-06-21"
$ws.Range("AH4").Value = 2
$ws.Range("AR4").Value = 18
$ws.Range("AV4").Value = 27
$ws.Range("BF4").Value = "2009-06-21"
$ws.Range("D5").Value = 82
$ws.Range("E5").Value = 41
$ws.Range("G5").Value = 0.5
$ws.Range("I5").Value = 38.1
$ws.Range("J5").Value = 83.5
$ws.Range("K5").Value = 0.457
$ws.Range("O5").Value = 19.9
$ws.Range("P5").Value = 25
$ws.Range("Q5").Value = 0.796
$ws.Range("U5").Value = 21.1
$ws.Range("Z5").Value = 20.8
$ws.Range("AA5").Value = 20.8
$ws.Range("AB5").Value = 102.2
$ws.Range("AC5").Value = -0.3
$ws.Range("AD5").Value = 1
$ws.Range("AE5").Value = 15
$ws.Range("AG5").Value = 15
$ws.Range("AI5").Value = 8
$ws.Range("AK5").Value = 15
$ws.Range("AO5").Value = 7
$ws.Range("AV5").Value = 20
$ws.Range("AW5").Value = 9
$ws.Range("BF5").Value = "2009-06-21"
$ws.Range("D6").Value = 82
$ws.Range("F6").Value = 16
$ws.Range("G6").Value = 0.805
$ws.Range("I6").Value = 36.9
$ws.Range("J6").Value = 78.7
$ws.Range("K6").Value = 0.468
$ws.Range("L6").Value = 8
$ws.Range("N6").Value = 0.393
$ws.Range("O6").Value = 18.6
$ws.Range("P6").Value = 24.5
$ws.Range("V6").Value = 12.7
$ws.Range("Y6").Value = 4.1
$ws.Range("Z6").Value = 20.3
$ws.Range("AA6").Value = 20.3
$ws.Range("AB6").Value = 100.3
$ws.Range("AC6").Value = 8.9
$ws.Range("AD6").Value = 1
$ws.Range("AO6").Value = 19
$ws.Range("AV6").Value = 6
$ws.Range("AW6").Value = 15
$ws.Range("AY6").Value = 6
$ws.Range("BF6").Value = "2009-06-21"
$ws.Range("AI7").Value = 7
$ws.Range("AV7").Value = 5
$ws.Range("BF7").Value = "2009-06-21"
$ws.Range("AF8").Value = 5
$ws.Range("AG8").Value = 5
$ws.Range("BF8").Value = "2009-06-21"
$ws.Range("AY9").Value = 7
$ws.Range("BF9").Value = "2009-06-21"
$ws.Range("BF10").Value = "2009-06-21"
$ws.Range("AO11").Value = 17
$ws.Range("BF11").Value = "2009-06-21"
$ws.Range("AH12").Value = 11
$ws.Range("AO12").Value = 18
$ws.Range("AV12").Value = 19
$ws.Range("AZ12").Value = 27
$ws.Range("BF12").Value = "2009-06-21"
$ws.Range("AH13").Value = 7
$ws.Range("BF13").Value = "2009-06-21"
$ws.Range("D14").Value = 82
$ws.Range("E14").Value = 65
$ws.Range("G14").Value = 0.793
$ws.Range("I14").Value = 40.3
$ws.Range("J14").Value = 85.09999999999999
$ws.Range("M14").Value = 18.5
$ws.Range("N14").Value = 0.361
$ws.Range("O14").Value = 19.6
$ws.Range("Q14").Value = 0.77
$ws.Range("S14").Value = 31.5
$ws.Range("V14").Value = 13.5
$ws.Range("Z14").Value = 20.7
$ws.Range("AA14").Value = 22.1
$ws.Range("AB14").Value = 106.9
$ws.Range("AC14").Value = 7.7
$ws.Range("AD14").Value = 1
$ws.Range("AN14").Value = 19
$ws.Range("AO14").Value = 11
$ws.Range("AS14").Value = 6
$ws.Range("AX14").Value = 10
$ws.Range("AY14").Value = 14
$ws.Range("BA14").Value = 8
$ws.Range("BC14").Value = 2
$ws.Range("BF14").Value = "2009-06-21"
$ws.Range("AN15").Value = 21
$ws.Range("BF15").Value = "2009-06-21"
$ws.Range("AK16").Value = 16
$ws.Range("BF16").Value = "2009-06-21"
$ws.Range("D17").Value = 82
$ws.Range("F17").Value = 48
$ws.Range("G17").Value = 0.415
$ws.Range("N17").Value = 0.363
$ws.Range("P17").Value = 25.2
$ws.Range("Q17").Value = 0.78
$ws.Range("S17").Value = 28.9
$ws.Range("T17").Value = 40.7
$ws.Range("V17").Value = 14.1
$ws.Range("W17").Value = 7.4
$ws.Range("Z17").Value = 24.2
$ws.Range("AC17").Value = -1.1
$ws.Range("AD17").Value = 1
$ws.Range("AF17").Value = 20
$ws.Range("AK17").Value = 27
$ws.Range("AP17").Value = 13
$ws.Range("AW17").Value = 12
$ws.Range("BF17").Value = "2009-06-21"
$ws.Range("BF18").Value = "2009-06-21"
$ws.Range("AF19").Value = 20
$ws.Range("AG19").Value = 20
$ws.Range("AH19").Value = 11
$ws.Range("AX19").Value = 15
$ws.Range("AY19").Value = 15
$ws.Range("BF19").Value = "2009-06-21"
$ws.Range("BF20").Value = "2009-06-21"
$ws.Range("AK21").Value = 28
$ws.Range("AN21").Value = 20
$ws.Range("AT21").Value = 8
$ws.Range("AW21").Value = 11
$ws.Range("AZ21").Value = 13
$ws.Range("BF21").Value = "2009-06-21"
$ws.Range("AO22").Value = 8
$ws.Range("AU22").Value = 20
$ws.Range("BF22").Value = "2009-06-21"
$ws.Range("AK23").Value = 17
$ws.Range("AO23").Value = 10
$ws.Range("BB23").Value = 10
$ws.Range("BF23").Value = "2009-06-21"
$ws.Range("BF24").Value = "2009-06-21"
$ws.Range("BF25").Value = "2009-06-21"
$ws.Range("D26").Value = 82
$ws.Range("F26").Value = 28
$ws.Range("G26").Value = 0.659
$ws.Range("I26").Value = 36.8
$ws.Range("J26").Value = 79.2
$ws.Range("K26").Value = 0.465
$ws.Range("L26").Value = 7.3
$ws.Range("M26").Value = 19
$ws.Range("O26").Value = 18.5
$ws.Range("P26").Value = 24.2
$ws.Range("Q26").Value = 0.765
$ws.Range("AA26").Value = 21.1
$ws.Range("AB26").Value = 99.40000000000001
$ws.Range("AC26").Value = 5.3
$ws.Range("AD26").Value = 1
$ws.Range("AK26").Value = 8
$ws.Range("AO26").Value = 20
$ws.Range("AQ26").Value = 17
$ws.Range("AT26").Value = 12
$ws.Range("AZ26").Value = 12
$ws.Range("BA26").Value = 12
$ws.Range("BF26").Value = "2009-06-21"
$ws.Range("BA27").Value = 11
$ws.Range("BF27").Value = "2009-06-21"
$ws.Range("AF28").Value = 5
$ws.Range("AG28").Value = 5
$ws.Range("AK28").Value = 7
$ws.Range("BF28").Value = "2009-06-21"
$ws.Range("AX29").Value = 16
$ws.Range("BF29").Value = "2009-06-21"
$ws.Range("AH30").Value = 7
$ws.Range("BF30").Value = "2009-06-21"
$ws.Range("AQ31").Value = 16
$ws.Range("AW31").Value = 8
$ws.Range("BF31").Value = "2009-06-21"
